# Append the new "2020-07-01" log entry as row 32 (SSA raw/clean data for July 1st).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Force the date-label column to store a literal text value (matching the
# "yyyy-mm-dd" text entries used for every other row) instead of letting
# Excel auto-convert it to a date serial number; ClearFormats() afterwards
# drops the temporary text number-format so the cell keeps the sheet's
# default (unstyled) appearance, same as the existing data rows.
$ws.Range("A32").NumberFormat = "@"
$ws.Range("A32").Value = "2020-07-01"
$ws.Range("A32").ClearFormats()

$ws.Range("B32").Value = 231770
$ws.Range("C32").Value = 289142
$ws.Range("D32").Value = 75005
$ws.Range("E32").Value = 28510
$ws.Range("F32").Value = 30.73
